$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.141.16"
$ws.Range("E2").Value = "  +0.10%  "
# Row 3
$ws.Range("D3").Value = "1.924.95"
$ws.Range("E3").Value = "  +0.61%  "
# Row 4
$ws.Range("E4").Value = "  -0.21%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.13"
$ws.Range("E5").Value = "  -2.78%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.19%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5065"
$ws.Range("E7").Value = "  -2.64%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4041"
$ws.Range("E8").Value = "  -0.58%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08283"
$ws.Range("E9").Value = "  -2.49%  "
# Row 10
$ws.Range("E10").Value = "  -0.99%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.78"
$ws.Range("E11").Value = "  -2.59%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.88"
$ws.Range("E12").Value = "  +3.71%  "
# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.915.39"
$ws.Range("E13").Value = "  +0.91%  "
# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.426"
$ws.Range("E14").Value = "  -0.24%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.344"
$ws.Range("E15").Value = "  -0.34%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.24%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.86"
$ws.Range("E17").Value = "  -2.26%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001102"
$ws.Range("E18").Value = "  -1.02%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06488"
$ws.Range("E19").Value = "  -3.13%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.67"
$ws.Range("E20").Value = "  +1.42%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -0.17%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.993"
$ws.Range("E22").Value = "  -0.26%  "
# Row 23
$ws.Range("D23").Value = "30.177.64"
$ws.Range("E23").Value = "  +0.13%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.32"
$ws.Range("E24").Value = "  -0.14%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.193"
$ws.Range("E25").Value = "  -1.48%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "22.23"
$ws.Range("E26").Value = "  +5.33%  "
# Row 27
$ws.Range("D27").Value = "2.127.81"
$ws.Range("E27").Value = "  +0.34%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.05"
$ws.Range("E28").Value = "  -0.39%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.385"
$ws.Range("E29").Value = "  -0.85%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.62"
$ws.Range("E30").Value = "  +0.83%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.133"
$ws.Range("E31").Value = "  +3.72%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1046"
$ws.Range("E32").Value = "  -1.91%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.021"
$ws.Range("E33").Value = "  +0.39%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.792"
$ws.Range("E34").Value = "  +5.03%  "
# Row 35
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.452"
$ws.Range("E35").Value = "  +5.70%  "
# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02454"
$ws.Range("E36").Value = "  -1.47%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06475"
$ws.Range("E37").Value = "  -1.45%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2169"
$ws.Range("E38").Value = "  -1.90%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.864"
$ws.Range("E39").Value = "  +0.78%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.193"
$ws.Range("E40").Value = "  -2.76%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6414"
$ws.Range("E41").Value = "  -1.50%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.41"
$ws.Range("E42").Value = "  -4.04%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.220"
$ws.Range("E43").Value = "  -1.65%  "
# Row 44
$ws.Range("E44").Value = "  -0.11%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.36"
$ws.Range("E45").Value = "  +0.72%  "
# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6016"
$ws.Range("E46").Value = "  -2.04%  "
# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.170"
$ws.Range("E47").Value = "  +4.14%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.650"
$ws.Range("E48").Value = "  -2.57%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.99"
$ws.Range("E49").Value = "  +0.10%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.220"
$ws.Range("E50").Value = "  -1.56%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.41"
$ws.Range("E51").Value = "  +0.14%  "

Write-Host "Applied 50 row updates (109 cell changes) to cryptos sheet"
